$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (162) down into the new rows
$ws.Range("A162:W162").Copy()
$ws.Range("A163:W165").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 163
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = "24/05/2024"
$ws.Range("C163").Value = 147.155
$ws.Range("D163").Value = 146.7375
$ws.Range("E163").Value = 143
$ws.Range("F163").Value = 138.5
$ws.Range("G163").Value = 64.75
$ws.Range("H163").Value = 132.4142857142857
$ws.Range("I163").Value = 67.83799999999999
$ws.Range("J163").Value = 62
$ws.Range("K163").Value = 116.925
$ws.Range("L163").Value = 152.126170212766
$ws.Range("M163").Value = 120
$ws.Range("N163").Value = 181.5
$ws.Range("O163").Value = 181.5
$ws.Range("P163").Value = 177.15
$ws.Range("Q163").Value = 171.5
$ws.Range("R163").Value = 99
$ws.Range("S163").Value = 156
$ws.Range("T163").Value = 0.3192307692307693
$ws.Range("U163").Value = 64.48999999999999
$ws.Range("V163").Value = 118.25
$ws.Range("W163").Value = 64.48999999999999

# Row 164
$ws.Range("A164").Value = 162
$ws.Range("B164").Value = "27/05/2024"
$ws.Range("C164").Value = 151.2021621621622
$ws.Range("D164").Value = 149.0714285714286
$ws.Range("E164").Value = 143.6666666666667
$ws.Range("F164").Value = 142
$ws.Range("G164").Value = 64.75
$ws.Range("H164").Value = 147.4166666666667
$ws.Range("I164").Value = 67.83799999999999
$ws.Range("J164").Value = 62
$ws.Range("K164").Value = 131.7718421052631
$ws.Range("L164").Value = 165.5733783783784
$ws.Range("M164").Value = 120
$ws.Range("N164").Value = 189.75
$ws.Range("O164").Value = 183.6666666666667
$ws.Range("P164").Value = 179.5
$ws.Range("Q164").Value = 175.5
$ws.Range("R164").Value = 99
$ws.Range("S164").Value = 156
$ws.Range("T164").Value = 0.3192307692307693
$ws.Range("U164").Value = 64.48999999999999
$ws.Range("V164").Value = 118.25
$ws.Range("W164").Value = 64.48999999999999

# Row 165
$ws.Range("A165").Value = 163
$ws.Range("B165").Value = "28/05/2024"
$ws.Range("C165").Value = 157.23
$ws.Range("D165").Value = 152.3333333333333
$ws.Range("E165").Value = 146.5
$ws.Range("F165").Value = 142
$ws.Range("G165").Value = 64.75
$ws.Range("H165").Value = 161.8733333333333
$ws.Range("I165").Value = 67.83799999999999
$ws.Range("J165").Value = 62
$ws.Range("K165").Value = 146.482
$ws.Range("L165").Value = 179.05864
$ws.Range("M165").Value = 120
$ws.Range("N165").Value = 196.5
$ws.Range("O165").Value = 183.6666666666667
$ws.Range("P165").Value = 179.5
$ws.Range("Q165").Value = 175.5
$ws.Range("R165").Value = 99
$ws.Range("S165").Value = 201
$ws.Range("T165").Value = 0.3192307692307693
$ws.Range("U165").Value = 64.48999999999999
$ws.Range("V165").Value = 118.25
$ws.Range("W165").Value = 64.48999999999999
